$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-26: column B (id), column C (speaker_variant); column D (is_prefered) cleared
$data = @(
    @("#end",          "End"),
    @("#apo",          "Apo"),
    @("#kup",          "Kup"),
    @("#aur",          "Aur"),
    @("#nat",          "Nat"),
    @("#mer",          "Mer"),
    @("#iris",         "Iris"),
    @("#ven",          "Ven"),
    @("#pallas",       "Pallas"),
    @("#cep",          "Cep"),
    @("#val",          "Val"),
    @("#gan",          "Gan"),
    @("#dia",          "Dia"),
    @("#apol",         "Apol"),
    @("#kupido",       "Kupido"),
    @("#pal",          "Pal"),
    @("#nacht",        "Nacht"),
    @("#r.-v.-nim",    "R. v. Nim"),
    @("#apollo",       "Apollo"),
    @("#r.v.-hard",    "R.v. Hard"),
    @("#diana.-pallas","Diana. Pallas"),
    @("#hard",         "Hard"),
    @("#isris",        "Isris"),
    @("#nac",          "Nac"),
    @("#endimion",     "Endimion")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = ""
}
